# "Test + Yorum Satırlarının Eklenmesi"
# Clear out the placeholder test values that were filled into the
# evaluation grid (B2:F7) for each "Program Ciktilari" row, leaving the
# per-row average formulas in column G untouched so they now compute
# against the (currently) empty inputs. Row heights are then left to
# auto-size to their (shorter) content, and the sheet's print setup is
# confirmed to fit one page wide/tall.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:F7").ClearContents()

$ws.Rows("1:7").AutoFit()

$ws.PageSetup.Zoom = 100
$ws.PageSetup.FitToPagesWide = 1
$ws.PageSetup.FitToPagesTall = 1
